$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$src = $ws.Range("F4:N5")
$dst = $ws.Range("G4:O5")
$src.Copy($dst)
